# Weekly update: insert two new daily price records into the
# "Hortaliza, Terminal Hortofrutícola Agro Chillán - Repollo" sheet.
#
# The sheet is a flat, date-ordered-by-insertion log where every row shares
# the same Mercado/Region/Categoria (columns A,B,C,E,F,G,N,Q,R); only the
# date (D), variedad (H), calidad (I), volumen (J), precios (K,L,M,P) and
# origen (O) vary row to row. Two brand-new rows are spliced into the
# middle of the table (not appended at the bottom), pushing every row
# below each insertion point down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert #1: new row at row 100 -----------------------------------
$ws.Rows.Item(100).Insert()

$ws.Range("A100").Value = 7
$ws.Range("B100").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C100").Value = "Ñuble"
$ws.Range("D100").Value = 44664
$ws.Range("E100").Value = 16
$ws.Range("F100").Value = 100112006
$ws.Range("G100").Value = "Repollo"
$ws.Range("H100").Value = "Crespo record"
$ws.Range("I100").Value = "Primera"
$ws.Range("J100").Value = 200
$ws.Range("K100").Value = 900
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = 950
$ws.Range("N100").Value = "`$/unidad"
$ws.Range("O100").Value = "Provincia de Diguillín"
$ws.Range("P100").Value = 950
$ws.Range("Q100").Value = 1
$ws.Range("R100").Value = "Hortaliza"

# --- Insert #2: new row at row 167 (post shift-#1 numbering) ---------
$ws.Rows.Item(167).Insert()

$ws.Range("A167").Value = 7
$ws.Range("B167").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C167").Value = "Ñuble"
$ws.Range("D167").Value = 44663
$ws.Range("E167").Value = 16
$ws.Range("F167").Value = 100112006
$ws.Range("G167").Value = "Repollo"
$ws.Range("H167").Value = "Crespo record"
$ws.Range("I167").Value = "Primera"
$ws.Range("J167").Value = 160
$ws.Range("K167").Value = 900
$ws.Range("L167").Value = 1000
$ws.Range("M167").Value = 950
$ws.Range("N167").Value = "`$/unidad"
$ws.Range("O167").Value = "Provincia de Diguillín"
$ws.Range("P167").Value = 950
$ws.Range("Q167").Value = 1
$ws.Range("R167").Value = "Hortaliza"
